# Corrects the Authors / Year / No. Primary Studies values that had been
# entered on the wrong rows: the "Wei et al." record (Year 2023, N=83)
# belongs with the row-3 title and the "Castaño-Villa et al." record
# (Year 2019, N=37) belongs with the row-2 title. Also drops the stray
# quote-prefix/fill formatting that had been picked up on I3, and moves
# the selection to where the editor left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch cell well outside the used range, used as a temporary holding
# spot so the B/C/E values (and B's cell formatting) can be swapped between
# row 2 and row 3 without clobbering each other.
$scratch = $ws.Cells.Item(100, 26)

# --- Authors (column B): swap value AND cell formatting ---
$ws.Cells.Item(2, 2).Cut($scratch)
$ws.Cells.Item(3, 2).Cut($ws.Cells.Item(2, 2))
$scratch.Cut($ws.Cells.Item(3, 2))

# --- Year (column C): swap values ---
$ws.Cells.Item(2, 3).Cut($scratch)
$ws.Cells.Item(3, 3).Cut($ws.Cells.Item(2, 3))
$scratch.Cut($ws.Cells.Item(3, 3))

# --- No. Primary Studies (column E): swap values ---
$ws.Cells.Item(2, 5).Cut($scratch)
$ws.Cells.Item(3, 5).Cut($ws.Cells.Item(2, 5))
$scratch.Cut($ws.Cells.Item(3, 5))

# Tidy up the scratch cell so it doesn't leave any trace behind.
$scratch.ClearFormats()
$scratch.ClearContents()

# I3 had picked up quote-prefix/fill formatting; strip it back to the
# default (the text value itself is unchanged).
$ws.Cells.Item(3, 9).ClearFormats()

# Update the on-screen selection / scroll position to match where the
# editor ended up.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("F13").Select()
